$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(24, 11).Value = -7.950704680006522
$ws.Cells.Item(25, 10).Value = -7.934690310484686
$ws.Cells.Item(25, 11).Value = 0.2049867170129686
$ws.Cells.Item(26, 9).Value = -8.018018899399499
$ws.Cells.Item(26, 10).Value = 0.1216581280981555
$ws.Cells.Item(26, 11).Value = 2.269229131468145
$ws.Cells.Item(27, 8).Value = -8.070357068665693
$ws.Cells.Item(27, 9).Value = 0.06931995883196146
$ws.Cells.Item(27, 10).Value = 2.216890962201951
$ws.Cells.Item(27, 11).Value = -1.223143443997281
$ws.Cells.Item(28, 7).Value = -8.04099600001669
$ws.Cells.Item(28, 8).Value = 0.0986810274809648
$ws.Cells.Item(28, 9).Value = 2.246252030850955
$ws.Cells.Item(28, 10).Value = -1.193782375348277
$ws.Cells.Item(28, 11).Value = -1.254965937862678
$ws.Cells.Item(29, 6).Value = -8.099654150585518
$ws.Cells.Item(29, 7).Value = 0.04002287691213546
$ws.Cells.Item(29, 8).Value = 2.187593880282126
$ws.Cells.Item(29, 9).Value = -1.252440525917107
$ws.Cells.Item(29, 10).Value = -1.313624088431507
$ws.Cells.Item(29, 11).Value = 0.7278401608117974
$ws.Cells.Item(30, 5).Value = -8.189565571715672
$ws.Cells.Item(30, 6).Value = -0.04988854421801681
$ws.Cells.Item(30, 7).Value = 2.097682459151973
$ws.Cells.Item(30, 8).Value = -1.342351947047259
$ws.Cells.Item(30, 9).Value = -1.40353550956166
$ws.Cells.Item(30, 10).Value = 0.6379287396816451
$ws.Cells.Item(30, 11).Value = 0.05533780647090653
$ws.Cells.Item(31, 4).Value = -8.238155430645202
$ws.Cells.Item(31, 5).Value = -0.09847840314754713
$ws.Cells.Item(31, 6).Value = 2.049092600222443
$ws.Cells.Item(31, 7).Value = -1.390941805976789
$ws.Cells.Item(31, 8).Value = -1.45212536849119
$ws.Cells.Item(31, 9).Value = 0.5893388807521147
$ws.Cells.Item(31, 10).Value = 0.006747947541376198
$ws.Cells.Item(31, 11).Value = 0.1024251963979083
$ws.Cells.Item(32, 3).Value = -8.527450494574166
$ws.Cells.Item(32, 4).Value = -0.387773467076512
$ws.Cells.Item(32, 5).Value = 1.759797536293478
$ws.Cells.Item(32, 6).Value = -1.680236869905754
$ws.Cells.Item(32, 7).Value = -1.741420432420155
$ws.Cells.Item(32, 8).Value = 0.3000438168231499
$ws.Cells.Item(32, 9).Value = -0.2825471163875887
$ws.Cells.Item(32, 10).Value = -0.1868698675310566
$ws.Cells.Item(32, 11).Value = -0.4472918755920416
$ws.Cells.Item(33, 2).Value = -9.539112594293229
$ws.Cells.Item(33, 3).Value = -1.399435566795574
$ws.Cells.Item(33, 4).Value = 0.7481354365744159
$ws.Cells.Item(33, 5).Value = -2.691898969624816
$ws.Cells.Item(33, 6).Value = -2.753082532139217
$ws.Cells.Item(33, 7).Value = -0.7116182828959121
$ws.Cells.Item(33, 8).Value = -1.294209216106651
$ws.Cells.Item(33, 9).Value = -1.198531967250119
$ws.Cells.Item(33, 10).Value = -1.458953975311104
$ws.Cells.Item(33, 11).Value = -0.6234959894775609
$ws.Cells.Item(34, 2).Value = -0.5626908242059827
$ws.Cells.Item(34, 3).Value = 1.584880179164007
$ws.Cells.Item(34, 4).Value = -1.855154227035225
$ws.Cells.Item(34, 5).Value = -1.916337789549626
$ws.Cells.Item(34, 6).Value = 0.1251264596936792
$ws.Cells.Item(34, 7).Value = -0.4574644735170594
$ws.Cells.Item(34, 8).Value = -0.3617872246605273
$ws.Cells.Item(34, 9).Value = -0.6222092327215123
$ws.Cells.Item(34, 10).Value = 0.2132487531120304
$ws.Cells.Item(34, 11).Value = 0.1132259840148222
$ws.Cells.Item(35, 2).Value = 1.582283197889339
$ws.Cells.Item(35, 3).Value = -1.857751208309893
$ws.Cells.Item(35, 4).Value = -1.918934770824294
$ws.Cells.Item(35, 5).Value = 0.1225294784190111
$ws.Cells.Item(35, 6).Value = -0.4600614547917275
$ws.Cells.Item(35, 7).Value = -0.3643842059351954
$ws.Cells.Item(35, 8).Value = -0.6248062139961804
$ws.Cells.Item(35, 9).Value = 0.2106517718373623
$ws.Cells.Item(35, 10).Value = 0.1106290027401541
$ws.Cells.Item(35, 11).Value = -0.2432290572937746
$ws.Cells.Item(36, 2).Value = -1.988887556949398
$ws.Cells.Item(36, 3).Value = -2.050071119463799
$ws.Cells.Item(36, 4).Value = -0.008606870220494112
$ws.Cells.Item(36, 5).Value = -0.5911978034312326
$ws.Cells.Item(36, 6).Value = -0.4955205545747006
$ws.Cells.Item(36, 7).Value = -0.7559425626356856
$ws.Cells.Item(36, 8).Value = 0.07951542319785709
$ws.Cells.Item(36, 9).Value = -0.0205073458993511
$ws.Cells.Item(36, 10).Value = -0.3743654059332798
$ws.Cells.Item(36, 11).Value = -0.2261723364979147
$ws.Cells.Item(37, 2).Value = -1.927016308533492
$ws.Cells.Item(37, 3).Value = 0.1144479407098132
$ws.Cells.Item(37, 4).Value = -0.4681429925009254
$ws.Cells.Item(37, 5).Value = -0.3724657436443933
$ws.Cells.Item(37, 6).Value = -0.6328877517053784
$ws.Cells.Item(37, 7).Value = 0.2025702341281644
$ws.Cells.Item(37, 8).Value = 0.1025474650309562
$ws.Cells.Item(37, 9).Value = -0.2513105950029725
$ws.Cells.Item(37, 10).Value = -0.1031175255676074
$ws.Cells.Item(37, 11).Value = 0.05431613594065687
$ws.Cells.Item(38, 2).Value = 0.286130677946619
$ws.Cells.Item(38, 3).Value = -0.2964602552641196
$ws.Cells.Item(38, 4).Value = -0.2007830064075875
$ws.Cells.Item(38, 5).Value = -0.4612050144685725
$ws.Cells.Item(38, 6).Value = 0.3742529713649702
$ws.Cells.Item(38, 7).Value = 0.274230202267762
$ws.Cells.Item(38, 8).Value = -0.07962785776616668
$ws.Cells.Item(38, 9).Value = 0.06856521166919841
$ws.Cells.Item(38, 10).Value = 0.2259988731774627
$ws.Cells.Item(38, 11).Value = -0.2746572130988596
$ws.Cells.Item(39, 2).Value = -0.2799830286212046
$ws.Cells.Item(39, 3).Value = -0.1843057797646726
$ws.Cells.Item(39, 4).Value = -0.4447277878256575
$ws.Cells.Item(39, 5).Value = 0.3907301980078852
$ws.Cells.Item(39, 6).Value = 0.290707428910677
$ws.Cells.Item(39, 7).Value = -0.06315063112325173
$ws.Cells.Item(39, 8).Value = 0.08504243831211336
$ws.Cells.Item(39, 9).Value = 0.2424760998203776
$ws.Cells.Item(39, 10).Value = -0.2581799864559446
$ws.Cells.Item(39, 11).Value = 0.02217630081466202
$ws.Cells.Item(40, 2).Value = -0.1907655701645288
$ws.Cells.Item(40, 3).Value = -0.4511875782255138
$ws.Cells.Item(40, 4).Value = 0.3842704076080289
$ws.Cells.Item(40, 5).Value = 0.2842476385108207
$ws.Cells.Item(40, 6).Value = -0.06961042152310797
$ws.Cells.Item(40, 7).Value = 0.07858264791225712
$ws.Cells.Item(40, 8).Value = 0.2360163094205214
$ws.Cells.Item(40, 9).Value = -0.2646397768558009
$ws.Cells.Item(40, 10).Value = 0.01571651041480578
$ws.Cells.Item(40, 11).Value = -0.1562496011445794
$ws.Cells.Item(41, 2).Value = -0.3709431027150346
$ws.Cells.Item(41, 3).Value = 0.4645148831185081
$ws.Cells.Item(41, 4).Value = 0.3644921140212999
$ws.Cells.Item(41, 5).Value = 0.01063405398737125
$ws.Cells.Item(41, 6).Value = 0.1588271234227363
$ws.Cells.Item(41, 7).Value = 0.3162607849310006
$ws.Cells.Item(41, 8).Value = -0.1843953013453217
$ws.Cells.Item(41, 9).Value = 0.09596098592528499
$ws.Cells.Item(41, 10).Value = -0.07600512563410017
$ws.Cells.Item(41, 11).Value = 0.2297419926957899
$ws.Cells.Item(42, 2).Value = 0.7003399477485883
$ws.Cells.Item(42, 3).Value = 0.6003171786513801
$ws.Cells.Item(42, 4).Value = 0.2464591186174515
$ws.Cells.Item(42, 5).Value = 0.3946521880528166
$ws.Cells.Item(42, 6).Value = 0.5520858495610809
$ws.Cells.Item(42, 7).Value = 0.0514297632847586
$ws.Cells.Item(42, 8).Value = 0.3317860505553653
$ws.Cells.Item(42, 9).Value = 0.1598199389959801
$ws.Cells.Item(42, 10).Value = 0.4655670573258702
$ws.Cells.Item(42, 11).Value = -0.1466231083065851
$ws.Cells.Item(43, 2).Value = 1.534006967818713
$ws.Cells.Item(43, 3).Value = 1.180148907784784
$ws.Cells.Item(43, 4).Value = 1.328341977220149
$ws.Cells.Item(43, 5).Value = 1.485775638728414
$ws.Cells.Item(43, 6).Value = 0.9851195524520913
$ws.Cells.Item(43, 7).Value = 1.265475839722698
$ws.Cells.Item(43, 8).Value = 1.093509728163313
$ws.Cells.Item(43, 9).Value = 1.399256846493203
$ws.Cells.Item(43, 10).Value = 0.7870666808607476
$ws.Cells.Item(43, 11).Value = 1.475402913200228
$ws.Cells.Item(44, 2).Value = 0.2665290763311319
$ws.Cells.Item(44, 3).Value = 0.414722145766497
$ws.Cells.Item(44, 4).Value = 0.5721558072747612
$ws.Cells.Item(44, 5).Value = 0.07149972099843901
$ws.Cells.Item(44, 6).Value = 0.3518560082690457
$ws.Cells.Item(44, 7).Value = 0.1798898967096605
$ws.Cells.Item(44, 8).Value = 0.4856370150395506
$ws.Cells.Item(44, 9).Value = -0.1265531505929047
$ws.Cells.Item(44, 10).Value = 0.5617830817465759
$ws.Cells.Item(45, 2).Value = 0.3874065102046941
$ws.Cells.Item(45, 3).Value = 0.5448401717129584
$ws.Cells.Item(45, 4).Value = 0.04418408543663613
$ws.Cells.Item(45, 5).Value = 0.3245403727072428
$ws.Cells.Item(45, 6).Value = 0.1525742611478576
$ws.Cells.Item(45, 7).Value = 0.4583213794777478
$ws.Cells.Item(45, 8).Value = -0.1538687861547076
$ws.Cells.Item(45, 9).Value = 0.534467446184773
$ws.Cells.Item(46, 2).Value = 0.7216226592872005
$ws.Cells.Item(46, 3).Value = 0.2209665730108782
$ws.Cells.Item(46, 4).Value = 0.5013228602814849
$ws.Cells.Item(46, 5).Value = 0.3293567487220997
$ws.Cells.Item(46, 6).Value = 0.6351038670519898
$ws.Cells.Item(46, 7).Value = 0.02291370141953453
$ws.Cells.Item(46, 8).Value = 0.7112499337590151
$ws.Cells.Item(47, 2).Value = -0.0137716841008341
$ws.Cells.Item(47, 3).Value = 0.2665846031697726
$ws.Cells.Item(47, 4).Value = 0.09461849161038739
$ws.Cells.Item(47, 5).Value = 0.4003656099402775
$ws.Cells.Item(47, 6).Value = -0.2118245556921778
$ws.Cells.Item(47, 7).Value = 0.4765116766473028
$ws.Cells.Item(48, 2).Value = 0.1211178113396557
$ws.Cells.Item(48, 3).Value = -0.0508483002197295
$ws.Cells.Item(48, 4).Value = 0.2548988181101606
$ws.Cells.Item(48, 5).Value = -0.3572913475222947
$ws.Cells.Item(48, 6).Value = 0.3310448848171859
$ws.Cells.Item(49, 2).Value = -0.02122337563261142
$ws.Cells.Item(49, 3).Value = 0.2845237426972787
$ws.Cells.Item(49, 4).Value = -0.3276664229351766
$ws.Cells.Item(49, 5).Value = 0.360669809404304
$ws.Cells.Item(50, 2).Value = 0.1707505866712637
$ws.Cells.Item(50, 3).Value = -0.4414395789611916
$ws.Cells.Item(50, 4).Value = 0.246896653378289
$ws.Cells.Item(51, 2).Value = -0.4496452844522253
$ws.Cells.Item(51, 3).Value = 0.2386909478872553
$ws.Cells.Item(52, 2).Value = 0.232387863851956
